# Auto-generated edit script: updates market price / profit columns (H-N)
# on specific rows across multiple sheets, per the target diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 100
$ws.Range("I5").Value = 69.44444
$ws.Range("J5").Value = 155
$ws.Range("K5").Value = 69.44444
$ws.Range("L5").Value = 155
$ws.Range("M5").Value = 45.55556
$ws.Range("N5").Value = -385

$ws.Range("H33").Value = 271.63635
$ws.Range("I33").Value = 185.5
$ws.Range("K33").Value = 185.5
$ws.Range("M33").Value = 43.5

$ws.Range("H40").Value = 2329.8845
$ws.Range("J40").Value = 2344.8333
$ws.Range("L40").Value = 2344.8333
$ws.Range("N40").Value = -2694.8333

$ws.Range("H51").Value = 38702.934
$ws.Range("J51").Value = 43852.23
$ws.Range("L51").Value = 43852.23
$ws.Range("N51").Value = -44820.23

$ws.Range("H116").Value = 2225901.8
$ws.Range("I116").Value = 11111111
$ws.Range("J116").Value = 4599.5
$ws.Range("K116").Value = 11111111
$ws.Range("L116").Value = 4599.5
$ws.Range("M116").Value = -11107669
$ws.Range("N116").Value = -11483.5

$ws.Range("H125").Value = 4178.125
$ws.Range("J125").Value = 4685.143
$ws.Range("L125").Value = 42166.287
$ws.Range("N125").Value = -47086.287


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 50145484
$ws.Range("J74").Value = 3503.25
$ws.Range("L74").Value = 3503.25
$ws.Range("N74").Value = -5251.25

$ws.Range("H77").Value = 50145484
$ws.Range("J77").Value = 3503.25
$ws.Range("L77").Value = 17516.25
$ws.Range("N77").Value = -26252.25


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 72996.39999999999
$ws.Range("J35").Value = 72996.39999999999
$ws.Range("L35").Value = 72996.39999999999
$ws.Range("N35").Value = -73616.39999999999

$ws.Range("H82").Value = 52109.875
$ws.Range("J82").Value = 104626.664
$ws.Range("L82").Value = 104626.664
$ws.Range("N82").Value = -105392.664

$ws.Range("H85").Value = 52109.875
$ws.Range("J85").Value = 104626.664
$ws.Range("L85").Value = 104626.664
$ws.Range("N85").Value = -107278.664

$ws.Range("H99").Value = 13897.794
$ws.Range("I99").Value = 15996.56
$ws.Range("K99").Value = 15996.56
$ws.Range("M99").Value = -14498.56

$ws.Range("H105").Value = 127549.78
$ws.Range("I105").Value = 278243.5
$ws.Range("K105").Value = 278243.5
$ws.Range("M105").Value = -276496.5

$ws.Range("H134").Value = 11513.3125
$ws.Range("I134").Value = 12283.071
$ws.Range("K134").Value = 36849.213
$ws.Range("M134").Value = -34314.213


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2769.95
$ws.Range("I31").Value = 1486.2
$ws.Range("J31").Value = 2953.3428
$ws.Range("K31").Value = 1486.2
$ws.Range("L31").Value = 2953.3428
$ws.Range("M31").Value = -1191.2
$ws.Range("N31").Value = -3543.3428

$ws.Range("H34").Value = 2769.95
$ws.Range("I34").Value = 1486.2
$ws.Range("J34").Value = 2953.3428
$ws.Range("K34").Value = 1486.2
$ws.Range("L34").Value = 2953.3428
$ws.Range("M34").Value = -1284.2
$ws.Range("N34").Value = -3357.3428


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 17243880
$ws.Range("I68").Value = 1965.6666
$ws.Range("K68").Value = 5896.9998
$ws.Range("M68").Value = -5085.9998

$ws.Range("H71").Value = 17243880
$ws.Range("I71").Value = 1965.6666
$ws.Range("K71").Value = 17690.9994
$ws.Range("M71").Value = -13634.9994

$ws.Range("H113").Value = 1055.2
$ws.Range("I113").Value = 555.61536
$ws.Range("J113").Value = 1350.409
$ws.Range("K113").Value = 1666.84608
$ws.Range("L113").Value = 4051.227
$ws.Range("M113").Value = 503.15392
$ws.Range("N113").Value = -8391.227000000001

$ws.Range("H134").Value = 2548.3076
$ws.Range("I134").Value = 2094
$ws.Range("K134").Value = 6282
$ws.Range("M134").Value = -1212


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 83.79310599999999
$ws.Range("I2").Value = 80.695656
$ws.Range("J2").Value = 95.666664
$ws.Range("K2").Value = 80.695656
$ws.Range("L2").Value = 95.666664
$ws.Range("M2").Value = 32.304344
$ws.Range("N2").Value = -321.666664

$ws.Range("H97").Value = 10558.091
$ws.Range("I97").Value = 13511.5
$ws.Range("J97").Value = 2682.3333
$ws.Range("K97").Value = 13511.5
$ws.Range("L97").Value = 2682.3333
$ws.Range("M97").Value = -13015.5
$ws.Range("N97").Value = -3674.3333

$ws.Range("H122").Value = 32300
$ws.Range("I122").Value = 49000
$ws.Range("K122").Value = 147000
$ws.Range("M122").Value = -144550


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 350000
$ws.Range("I40").Value = 350000
$ws.Range("K40").Value = 350000
$ws.Range("M40").Value = -349864

$ws.Range("H61").Value = 2359.45
$ws.Range("I61").Value = 1828
$ws.Range("J61").Value = 3009
$ws.Range("K61").Value = 1828
$ws.Range("L61").Value = 3009
$ws.Range("M61").Value = -1626
$ws.Range("N61").Value = -3413

$ws.Range("H68").Value = 4199.75
$ws.Range("I68").Value = 2744.2
$ws.Range("J68").Value = 4861.364
$ws.Range("K68").Value = 2744.2
$ws.Range("L68").Value = 4861.364
$ws.Range("M68").Value = -1995.2
$ws.Range("N68").Value = -6359.364

$ws.Range("H71").Value = 4199.75
$ws.Range("I71").Value = 2744.2
$ws.Range("J71").Value = 4861.364
$ws.Range("K71").Value = 13721
$ws.Range("L71").Value = 24306.82
$ws.Range("M71").Value = -9977
$ws.Range("N71").Value = -31794.82

$ws.Range("H93").Value = 3020.9375
$ws.Range("I93").Value = 3410.7
$ws.Range("J93").Value = 2371.3333
$ws.Range("K93").Value = 3410.7
$ws.Range("L93").Value = 2371.3333
$ws.Range("M93").Value = -2162.7
$ws.Range("N93").Value = -4867.3333

$ws.Range("H113").Value = 2359.45
$ws.Range("I113").Value = 1828
$ws.Range("J113").Value = 3009
$ws.Range("K113").Value = 1828
$ws.Range("L113").Value = 3009
$ws.Range("M113").Value = 342
$ws.Range("N113").Value = -7349

$ws.Range("H122").Value = 6150
$ws.Range("I122").Value = 5500
$ws.Range("J122").Value = 6475
$ws.Range("K122").Value = 16500
$ws.Range("L122").Value = 19425
$ws.Range("M122").Value = -14050
$ws.Range("N122").Value = -24325


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2791.6296
$ws.Range("I113").Value = 1195.9375
$ws.Range("K113").Value = 3587.8125
$ws.Range("M113").Value = -1417.8125

$ws.Range("H132").Value = 46870.43
$ws.Range("I132").Value = 67299.55499999999
$ws.Range("K132").Value = 201898.665
$ws.Range("M132").Value = -199368.665

$ws.Range("H136").Value = 3268.4866
$ws.Range("I136").Value = 2303.724
$ws.Range("K136").Value = 6911.172
$ws.Range("M136").Value = -4361.172

